$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 245 so the existing weekly records
# (currently in rows 245-249) shift down to rows 247-251, making room
# for the newest week's data (Especial + Segunda) at rows 245-246.
$ws.Rows.Item(245).Insert()
$ws.Rows.Item(245).Insert()

# --- New row 245: Femacal de La Calera, Frutilla, Especial (week 2022-02-03) ---
$ws.Cells.Item(245, 1).Value = 3
$ws.Cells.Item(245, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(245, 3).Value = "Coquimbo"
$ws.Cells.Item(245, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(245, 4).Value = 44595
$ws.Cells.Item(245, 5).Value = 5
$ws.Cells.Item(245, 6).Value = "Fruta"
$ws.Cells.Item(245, 7).Value = 100101
$ws.Cells.Item(245, 8).Value = "Berries"
$ws.Cells.Item(245, 9).Value = 100112025
$ws.Cells.Item(245, 10).Value = "Frutilla"
$ws.Cells.Item(245, 11).Value = "Sin especificar"
$ws.Cells.Item(245, 12).Value = "Especial"
$ws.Cells.Item(245, 13).Value = 65
$ws.Cells.Item(245, 14).Value = 6000
$ws.Cells.Item(245, 15).Value = 6000
$ws.Cells.Item(245, 16).Value = 6000
$ws.Cells.Item(245, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(245, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(245, 19).Value = 857
$ws.Cells.Item(245, 20).Value = 7

# --- New row 246: Femacal de La Calera, Frutilla, Segunda (week 2022-02-03) ---
$ws.Cells.Item(246, 1).Value = 3
$ws.Cells.Item(246, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(246, 3).Value = "Coquimbo"
$ws.Cells.Item(246, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(246, 4).Value = 44595
$ws.Cells.Item(246, 5).Value = 5
$ws.Cells.Item(246, 6).Value = "Fruta"
$ws.Cells.Item(246, 7).Value = 100101
$ws.Cells.Item(246, 8).Value = "Berries"
$ws.Cells.Item(246, 9).Value = 100112025
$ws.Cells.Item(246, 10).Value = "Frutilla"
$ws.Cells.Item(246, 11).Value = "Sin especificar"
$ws.Cells.Item(246, 12).Value = "Segunda"
$ws.Cells.Item(246, 13).Value = 58
$ws.Cells.Item(246, 14).Value = 4000
$ws.Cells.Item(246, 15).Value = 4000
$ws.Cells.Item(246, 16).Value = 4000
$ws.Cells.Item(246, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(246, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(246, 19).Value = 571
$ws.Cells.Item(246, 20).Value = 7

Write-Output ("Final UsedRange: " + $ws.UsedRange.Address())
